# The edit reshuffles the 29 data rows (rows 2-30) of the sheet: every row's
# full record (all columns A:R) moves to a different row position, while the
# set of records itself is unchanged. Build the mapping of new row -> old row
# and rewrite the data block accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# newRow -> oldRow
$map = @{2=21; 3=7; 4=22; 5=19; 6=27; 7=28; 8=8; 9=9; 10=29; 11=6; 12=5; 13=11; 14=24; 15=25; 16=15; 17=16; 18=30; 19=18; 20=3; 21=4; 22=20; 23=2; 24=26; 25=10; 26=23; 27=17; 28=12; 29=13; 30=14}

$firstRow = 2
$lastRow = 30
$firstCol = 1
$lastCol = 18

# Read the whole data block first so writes don't clobber values we still need to read.
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$original = $srcRange.Value2

$numRows = $lastRow - $firstRow + 1
$numCols = $lastCol - $firstCol + 1
$newArr = New-Object 'object[,]' $numRows, $numCols

for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $map[$newRow]
    $srcIdx = $oldRow - $firstRow + 1
    $dstIdx = $newRow - $firstRow + 1
    for ($c = 1; $c -le $numCols; $c++) {
        $newArr[$dstIdx - 1, $c - 1] = $original[$srcIdx, $c]
    }
}

$dstRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$dstRange.Value2 = $newArr
